{"js": "// Update the date line and every two-digit x two-digit multiplication\n// expression in the table to the new day / new set of problems.\nconst replacements = [\n  [\"2025-12-17 Wednesday\", \"2025-12-18 Thursday\"],\n  [\"58\u00d761=\", \"33\u00d725=\"],\n  [\"60\u00d739=\", \"73\u00d746=\"],\n  [\"46\u00d781=\", \"47\u00d757=\"],\n  [\"31\u00d736=\", \"84\u00d760=\"],\n  [\"48\u00d774=\", \"87\u00d756=\"],\n  [\"21\u00d794=\", \"85\u00d750=\"],\n  [\"77\u00d713=\", \"61\u00d736=\"],\n  [\"22\u00d758=\", \"37\u00d795=\"],\n  [\"83\u00d712=\", \"62\u00d797=\"],\n  [\"50\u00d746=\", \"45\u00d717=\"],\n  [\"38\u00d741=\", \"27\u00d777=\"],\n  [\"88\u00d743=\", \"20\u00d719=\"],\n  [\"35\u00d790=\", \"82\u00d774=\"],\n  [\"47\u00d799=\", \"46\u00d780=\"],\n  [\"61\u00d751=\", \"79\u00d732=\"],\n  [\"78\u00d783=\", \"18\u00d720=\"],\n  [\"32\u00d749=\", \"49\u00d780=\"],\n  [\"78\u00d732=\", \"33\u00d734=\"],\n  [\"66\u00d741=\", \"21\u00d741=\"],\n  [\"70\u00d788=\", \"36\u00d741=\"],\n  [\"25\u00d744=\", \"57\u00d764=\"],\n  [\"78\u00d784=\", \"64\u00d789=\"],\n  [\"31\u00d789=\", \"59\u00d752=\"],\n  [\"68\u00d756=\", \"20\u00d793=\"],\n  [\"39\u00d766=\", \"73\u00d768=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every two-digit x two-digit multiplication\n# expression in the table to the new day / new set of problems.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-17 Wednesday\", \"2025-12-18 Thursday\"),\n    @(\"58\u00d761=\", \"33\u00d725=\"),\n    @(\"60\u00d739=\", \"73\u00d746=\"),\n    @(\"46\u00d781=\", \"47\u00d757=\"),\n    @(\"31\u00d736=\", \"84\u00d760=\"),\n    @(\"48\u00d774=\", \"87\u00d756=\"),\n    @(\"21\u00d794=\", \"85\u00d750=\"),\n    @(\"77\u00d713=\", \"61\u00d736=\"),\n    @(\"22\u00d758=\", \"37\u00d795=\"),\n    @(\"83\u00d712=\", \"62\u00d797=\"),\n    @(\"50\u00d746=\", \"45\u00d717=\"),\n    @(\"38\u00d741=\", \"27\u00d777=\"),\n    @(\"88\u00d743=\", \"20\u00d719=\"),\n    @(\"35\u00d790=\", \"82\u00d774=\"),\n    @(\"47\u00d799=\", \"46\u00d780=\"),\n    @(\"61\u00d751=\", \"79\u00d732=\"),\n    @(\"78\u00d783=\", \"18\u00d720=\"),\n    @(\"32\u00d749=\", \"49\u00d780=\"),\n    @(\"78\u00d732=\", \"33\u00d734=\"),\n    @(\"66\u00d741=\", \"21\u00d741=\"),\n    @(\"70\u00d788=\", \"36\u00d741=\"),\n    @(\"25\u00d744=\", \"57\u00d764=\"),\n    @(\"78\u00d784=\", \"64\u00d789=\"),\n    @(\"31\u00d789=\", \"59\u00d752=\"),\n    @(\"68\u00d756=\", \"20\u00d793=\"),\n    @(\"39\u00d766=\", \"73\u00d768=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
